# commit: "fixed script for run NF"
#
# The Scene config's FilePath column pointed one directory level too high
# (../../NFDataCfg/Ini/Scene/N.xml). The NF server actually resolves these
# paths relative to one level up, so drop the extra "../" to fix the
# generated NF startup/run script.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

$ws.Range("F10").Value = "../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("F11").Value = "../NFDataCfg/Ini/Scene/2.xml"
$ws.Range("F12").Value = "../NFDataCfg/Ini/Scene/3.xml"
$ws.Range("F13").Value = "../NFDataCfg/Ini/Scene/4.xml"
$ws.Range("F14").Value = "../NFDataCfg/Ini/Scene/5.xml"
$ws.Range("F15").Value = "../NFDataCfg/Ini/Scene/6.xml"

# Leave the cursor where the author left it when they saved the file.
$ws.Range("F23").Select()
